# "Simple Page Object example"
#
# The source deck had an accidental duplicate of the "Locators in Selenium"
# picture on slide 6 (two <p:pic> shapes pointing at the same image,
# "Picture 6" and "Picture 8"). Remove the stray duplicate ("Picture 8"),
# keeping the original ("Picture 6") intact.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -eq "Picture 8") {
        $shape.Delete()
    }
}

# Incidental authoring artifact: opening/saving the deck in PowerPoint
# stamps an (empty) slide-guide list extension on the presentation.
# Touch the Guides collection so it gets recorded if the host supports it.
try {
    $p.Guides.Add(1, 3.0) | Out-Null
} catch {
}
